# Reformatting / audit pass over the GRNmap test workbook:
#  - the stray "Sheet" row (id/value = 3/4 placeholder leftover) in the
#    optimization_parameters sheet is removed, which also drops the now-unused
#    "Sheet" shared string and shifts everything below it up by one row.
#  - the active sheet moves from optimization_parameters to threshold_b.

$wb = $excel.ActiveWorkbook

$optParams = $wb.Worksheets.Item("optimization_parameters")
$optParams.Rows(16).Delete()

$thresholdB = $wb.Worksheets.Item("threshold_b")
$thresholdB.Activate()
